$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("G2").Value = 1.21
$ws.Range("H2").Value = 16
$ws.Range("I2").Value = 19.5
$ws.Range("U2").Value = 2
$ws.Range("W2").Value = 5.8

# --- Row 4 updates ---
$ws.Range("F4").Value = 2.52
$ws.Range("G4").Value = 2.54
$ws.Range("H4").Value = 3.4
$ws.Range("X4").Value = 10

# --- Row 5 updates ---
$ws.Range("AO5").Value = 29

# --- Row 6 updates ---
$ws.Range("F6").Value = 4.7
$ws.Range("G6").Value = 5.6
$ws.Range("H6").Value = 1.76
$ws.Range("I6").Value = 1.92
$ws.Range("J6").Value = 3.7
$ws.Range("N6").Value = 3.5
$ws.Range("P6").Value = 1.86
$ws.Range("Q6").Value = 1.93
$ws.Range("R6").Value = 1.33
$ws.Range("V6").Value = 2.08
$ws.Range("W6").Value = 1.22
$ws.Range("AD6").Value = 10.5
$ws.Range("AO6").Value = 13.5

# --- Row 7 updates ---
$ws.Range("S7").Value = 1.02

# --- Row 8 updates ---
$ws.Range("F8").Value = 1.44
$ws.Range("H8").Value = 9.8
$ws.Range("L8").Value = 1.41
$ws.Range("P8").Value = 1.86
$ws.Range("V8").Value = 1.1
$ws.Range("W8").Value = 3.15
$ws.Range("Y8").Value = 26
$ws.Range("AD8").Value = 40
$ws.Range("AL8").Value = 55
$ws.Range("AN8").Value = 8.6

# --- Row 9 updates ---
$ws.Range("F9").Value = 1.62
$ws.Range("G9").Value = 1.63
$ws.Range("L9").Value = 1.33
$ws.Range("O9").Value = 1.28
$ws.Range("R9").Value = 1.45
$ws.Range("U9").Value = 2.06
$ws.Range("V9").Value = 1.19
$ws.Range("W9").Value = 2.6
$ws.Range("X9").Value = 17.5
$ws.Range("Y9").Value = 21
$ws.Range("Z9").Value = 50
$ws.Range("AA9").Value = 180
$ws.Range("AB9").Value = 8.8
$ws.Range("AC9").Value = 9.4
$ws.Range("AD9").Value = 23
$ws.Range("AE9").Value = 85
$ws.Range("AF9").Value = 9.2
$ws.Range("AH9").Value = 21
$ws.Range("AI9").Value = 80
$ws.Range("AK9").Value = 15.5
$ws.Range("AL9").Value = 34
$ws.Range("AO9").Value = 95

# --- New Row 10 ---
$ws.Range("A10").Value = "Portuguese Segunda Liga"
# Force B10 to remain text "2026-01-08" instead of being auto-converted
# to a date serial number, then strip the temporary text number-format
# so the cell ends up unstyled (matching the other rows in this column).
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "2026-01-08"
$ws.Range("B10").ClearFormats()
$ws.Range("C10").Value = "17:15:00"
$ws.Range("D10").Value = "Lusitania Futebol Clube"
$ws.Range("E10").Value = "Leixoes"
$ws.Range("F10").Value = 1.98
$ws.Range("G10").Value = 2.62
$ws.Range("H10").Value = 3.35
$ws.Range("I10").Value = 5.2
$ws.Range("J10").Value = 2.84
$ws.Range("K10").Value = 5.3
$ws.Range("L10").Value = 1.01
$ws.Range("M10").Value = 1.01
$ws.Range("N10").Value = 1.53
$ws.Range("O10").Value = 1.01
$ws.Range("P10").Value = 1.53
$ws.Range("Q10").Value = 2.1
$ws.Range("R10").Value = 1.19
$ws.Range("S10").Value = 3.7
$ws.Range("T10").Value = 1.01
$ws.Range("U10").Value = 1.01
$ws.Range("V10").Value = 1.24
$ws.Range("W10").Value = 1.61
$ws.Range("X10").Value = 1000
$ws.Range("Y10").Value = 1000
$ws.Range("Z10").Value = 1000
$ws.Range("AA10").Value = 1000
$ws.Range("AB10").Value = 1000
$ws.Range("AC10").Value = 1000
$ws.Range("AD10").Value = 1000
$ws.Range("AE10").Value = 1000
$ws.Range("AF10").Value = 1000
$ws.Range("AG10").Value = 1000
$ws.Range("AH10").Value = 1000
$ws.Range("AI10").Value = 1000
$ws.Range("AJ10").Value = 1000
$ws.Range("AK10").Value = 1000
$ws.Range("AL10").Value = 1000
$ws.Range("AM10").Value = 1000
$ws.Range("AN10").Value = 1000
$ws.Range("AO10").Value = 1000

Write-Host "Edits applied successfully"
